$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each changed cell is re-written as literal text (apostrophe-prefixed to
# defeat the automatic number/date inference for numeric-looking strings
# like "1.00" or "43.485.24"), then the quote-prefix cell style that the
# apostrophe leaves behind is reset back to Normal so the saved file
# carries no extra formatting, matching the source diff exactly.

$ws.Range("D2").Value = "'43.485.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.257.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.44%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'231.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.643"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.33%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'64.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.45%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.41%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -7.68%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'56.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.58%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'26.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.49%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.588.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.80%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.254.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.345.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.48%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0967"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'73.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.64%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.67%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'247.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.50%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +10.45%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.13%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.09%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'173.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'21.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.60%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Kaspa"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'0.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'ImmutableX"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.27%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0679"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.62%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.89%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.19%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'6.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.40%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.55%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.65%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'8.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'FTXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'4.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.92%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'InjectiveProtocol"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'17.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.42%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Celestia"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'10.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.32%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'97.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.72%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0936"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.68%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -4.79%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.435.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.90%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.85%  "
$ws.Range("E51").Style = "Normal"
